# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "60.712.02"
$ws.Range("E2").Value  = "  -0.35%  "

$ws.Range("D3").Value  = "2.906.54"
$ws.Range("E3").Value  = "  -0.48%  "

$ws.Range("D4").Value  = "0.999"
$ws.Range("E4").Value  = "  -0.09%  "

$ws.Range("D5").Value  = "588.44"
$ws.Range("E5").Value  = "  +0.33%  "

$ws.Range("D6").Value  = "144.18"
$ws.Range("E6").Value  = "  -0.89%  "

$ws.Range("E7").Value  = "  +0.04%  "

$ws.Range("D8").Value  = "0.504"
$ws.Range("E8").Value  = "  -0.21%  "

$ws.Range("E9").Value  = "  +0.47%  "

$ws.Range("E10").Value = "  -2.58%  "

$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").Value = "33.34"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "3.385.99"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("D16").Value = "60.615.16"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "6.66"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").Value = "2.908.79"

$ws.Range("D19").Value = "430.37"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").Value = "13.31"
$ws.Range("E20").Value = "  -2.25%  "

$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").Value = "7.07"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").Value = "81.13"
$ws.Range("E23").Value = "  +0.53%  "

$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("E25").Value = "  -3.40%  "

$ws.Range("D26").Value = "11.74"
$ws.Range("E26").Value = "  -2.18%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +4.51%  "

$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("E30").Value = "  -3.60%  "

$ws.Range("D31").Value = "26.46"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "0.108"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("D34").Value = "0.0₃0854"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").Value = "5.60"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").Value = "2.96"
$ws.Range("E37").Value = "  -2.04%  "

$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("E39").Value = "  -5.18%  "

$ws.Range("D40").Value = "8.52"
$ws.Range("E40").Value = "  -1.60%  "

$ws.Range("D41").Value = "41.25"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "0.281"
$ws.Range("E42").Value = "  -5.94%  "

$ws.Range("D43").Value = "374.02"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "2.693.21"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -3.13%  "

$ws.Range("D46").Value = "133.79"
$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("D48").Value = "23.65"
$ws.Range("E48").Value = "  -3.46%  "

$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("E50").Value = "  -3.67%  "

$ws.Range("E51").Value = "  -1.16%  "
